# Append 25 new customer records (rows 277-301) to Sheet0,
# matching "R22 UAT2 - Regression" data fixture update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: TC, Customer_ID, PD, CD (optional), extra code (optional)
$newRows = @(
    @("118518", "17866752", "6020", "", ""),
    @("118498", "17866753", "1001", "", ""),
    @("118452", "17866754", "1001", "", ""),
    @("118518", "17866755", "1035", "", ""),
    @("118452", "17866756", "1150", "", ""),
    @("118448", "17866757", "1068", "", ""),
    @("118448", "17866759", "1005", "", ""),
    @("118518", "17866761", "6020", "", ""),
    @("118498", "17866762", "1001", "", ""),
    @("118452", "17866763", "1001", "", ""),
    @("118518", "17866764", "6005", "", ""),
    @("118448", "17866766", "1047", "", ""),
    @("118452", "17866767", "1150", "", ""),
    @("118448", "17866768", "1068", "", ""),
    @("118448", "17866770", "1005", "", ""),
    @("118463", "17868014", "1010", "CD Karobar Fin", "11"),
    @("118463", "17868016", "1007", "CD Karobar Fin", "11"),
    @("118463", "17868020", "1010", "CD Karobar Fin", "11"),
    @("118463", "17868021", "1003", "CD Karobar Fin", "11"),
    @("118463", "17868022", "1007", "CD Karobar Fin", "11"),
    @("118463", "17868023", "1010", "CD Karobar Fin", "11"),
    @("118463", "17868024", "1003", "CD Karobar Fin", "11"),
    @("118463", "17868025", "1007", "CD Karobar Fin", "11"),
    @("118463", "17868029", "1010", "CD F.E.25", "21"),
    @("118463", "17868030", "1003", "CD F.E.25", "21")
)

$startRow = 277
$endRow = $startRow + $newRows.Count - 1

# These columns hold ID-like values ("118518", "17866752", ...) that must
# stay text (as the rest of the sheet already does) instead of being
# auto-converted to numbers, so format them as Text before writing.
$ws.Range("A" + $startRow + ":C" + $endRow).NumberFormat = "@"
$ws.Range("E292:E301").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    if ($data[3] -ne "") {
        $ws.Cells.Item($r, 4).Value = $data[3]
    }
    if ($data[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = $data[4]
    }
}
